$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPt = 12700
$left   = 1069848 / $emuPerPt
$top    = 5863905 / $emuPerPt
$width  = 7891272 / $emuPerPt
$height = 489036 / $emuPerPt

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "Subtítulo 2"
$shp.TextFrame.TextRange.Text = "TCC – Arquitetura de Software Distribuído – PUC-MG"
$shp.TextFrame.TextRange.Font.Size = 11
